# Added column for "Status as of July 4, 2025" with a dropdown validation
# list backed by a new hidden sheet "DropdownOptions".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Create the hidden helper sheet (placed right after Sheet1) with the
#    dropdown option values.
$dropdown = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$dropdown.Name = "DropdownOptions"

$dropdown.Range("A1").Value = "0% - 10%"
$dropdown.Range("A2").Value = "11% - 25%"
$dropdown.Range("A3").Value = "26% - 50%"
$dropdown.Range("A4").Value = "51% - 75%"
$dropdown.Range("A5").Value = "76% - 90%"
$dropdown.Range("A6").Value = "91% - 99%"
# Leading apostrophe forces literal text so Excel doesn't auto-convert
# "100%" into the number 1 with a percentage format; reset the style
# back to Normal afterwards so the cell keeps the default (unstyled) look.
$dropdown.Range("A7").Value = "'100%"
$dropdown.Range("A7").Style = "Normal"

$dropdown.Visible = $false

# 2. Add the new header column on Sheet1.
$ws.Range("AH1").Value = "Status as of July 4, 2025"

# 3. Clean up the stray empty cells that were previously present (N2,
#    AE2, AF2 become fully blank again).
$ws.Range("N2").ClearContents()
$ws.Range("AE2").ClearContents()
$ws.Range("AF2").ClearContents()

# 4. Add a dropdown (data validation list) on AH2 pointing at
#    DropdownOptions!$A$1:$A$7
$target = $ws.Range("AH2")
$target.Validation.Add(3, 1, 1, "=DropdownOptions!`$A`$1:`$A`$7")
$target.Validation.IgnoreBlank = $true
$target.Validation.InCellDropdown = $true
$target.Validation.ShowInput = $false
$target.Validation.ShowError = $false

# 5. Keep Sheet1 as the active/selected sheet (creating the new sheet
#    shifts focus onto it by default).
$ws.Activate()
